$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.717.92'
$ws.Range('E2').Value = '  +0.13%  '

# Row 3
$ws.Range('D3').Value = '1.534.33'
$ws.Range('E3').Value = '  -1.45%  '

# Row 4
$ws.Range('E4').Value = '  -0.20%  '

# Row 5
$ws.Range('D5').Value = '''205.83'
$ws.Range('E5').Value = '  +0.04%  '

# Row 6
$ws.Range('E6').Value = '  -0.94%  '

# Row 8
$ws.Range('D8').Value = '''21.36'
$ws.Range('E8').Value = '  -2.64%  '

# Row 9
$ws.Range('E9').Value = '  -1.18%  '

# Row 10
$ws.Range('E10').Value = '  -0.53%  '

# Row 11
$ws.Range('E11').Value = '  -1.16%  '

# Row 12
$ws.Range('D12').Value = '1.752.35'
$ws.Range('E12').Value = '  -1.45%  '

# Row 13
$ws.Range('D13').Value = '1.532.38'
$ws.Range('E13').Value = '  -1.69%  '

# Row 14
$ws.Range('D14').Value = '''3.67'
$ws.Range('E14').Value = '  -1.57%  '

# Row 15
$ws.Range('D15').Value = '''0.506'
$ws.Range('E15').Value = '  -0.99%  '

# Row 16
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '26.712.70'
$ws.Range('E16').Value = '  -0.08%  '

# Row 17
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '''61.26'
$ws.Range('E17').Value = '  -0.53%  '

# Row 18
$ws.Range('D18').Value = '''212.21'
$ws.Range('E18').Value = '  -0.42%  '

# Row 19
$ws.Range('E19').Value = '  +1.15%  '

# Row 20
$ws.Range('E20').Value = '  -1.78%  '

# Row 21
$ws.Range('E21').Value = '  -0.12%  '

# Row 22
$ws.Range('D22').Value = '''4.00'
$ws.Range('E22').Value = '  -1.74%  '

# Row 23
$ws.Range('D23').Value = '''9.08'
$ws.Range('E23').Value = '  -2.80%  '

# Row 24
$ws.Range('D24').Value = '''1.95'
$ws.Range('E24').Value = '  -2.55%  '

# Row 25
$ws.Range('D25').Value = '''151.76'
$ws.Range('E25').Value = '  -0.29%  '

# Row 26
$ws.Range('E26').Value = '  -3.39%  '

# Row 27
$ws.Range('D27').Value = '''14.84'
$ws.Range('E27').Value = '  +0.20%  '

# Row 28
$ws.Range('E28').Value = '  -0.20%  '

# Row 29
$ws.Range('E29').Value = '  -1.08%  '

# Row 30
$ws.Range('E30').Value = '  -0.93%  '

# Row 31
$ws.Range('E31').Value = '  -1.97%  '

# Row 32
$ws.Range('E32').Value = '  +2.89%  '

# Row 33
$ws.Range('D33').Value = '1.357.27'
$ws.Range('E33').Value = '  -2.00%  '

# Row 34
$ws.Range('E34').Value = '  +0.14%  '

# Row 35
$ws.Range('E35').Value = '  -3.39%  '

# Row 36
$ws.Range('D36').Value = '''0.938'
$ws.Range('E36').Value = '  +0.80%  '

# Row 37
$ws.Range('E37').Value = '  -0.51%  '

# Row 38
$ws.Range('E38').Value = '  +0.39%  '

# Row 39
$ws.Range('D39').Value = '''0.521'
$ws.Range('E39').Value = '  +0.42%  '

# Row 40
$ws.Range('E40').Value = '  -1.65%  '

# Row 41
$ws.Range('E41').Value = '  -0.17%  '

# Row 42
$ws.Range('E42').Value = '  +5.41%  '

# Row 43
$ws.Range('E43').Value = '  -0.05%  '

# Row 44
$ws.Range('E44').Value = '  +0.43%  '

# Row 45
$ws.Range('E45').Value = '  -0.82%  '

# Row 46
$ws.Range('D46').Value = '''62.56'
$ws.Range('E46').Value = '  -0.77%  '

# Row 47
$ws.Range('D47').Value = '1.666.53'
$ws.Range('E47').Value = '  -1.47%  '

# Row 48
$ws.Range('D48').Value = '''85.41'
$ws.Range('E48').Value = '  -0.13%  '

# Row 49
$ws.Range('D49').Value = '''0.0506'
$ws.Range('E49').Value = '  +2.78%  '

# Row 50
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0968'
$ws.Range('E50').Value = '  -0.48%  '

# Row 51
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '''0.0943'
$ws.Range('E51').Value = '  -0.46%  '

